$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, $cellAddr, $text) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# --- Sheets ---
$wsLogin   = $wb.Worksheets.Item("loginQA")
$wsUnit    = $wb.Worksheets.Item("Unit to Reconcile")
$wsPO      = $wb.Worksheets.Item("PO_Detail")
$wsAccept  = $wb.Worksheets.Item("acceptQueue")
$wsOut     = $wb.Worksheets.Item("Unit_to_Reconcile_Output")

# --- Reconcile data update (10-May-2015 reconcile: new unit / PO info) ---
Set-TextValue $wsOut "A2" "00997989"
Set-TextValue $wsOut "B2" "00997989 - 2017 Ford F-450 Chassis XL 4x2 SD Crew Cab 179 in. WB DRW (W4G)"
Set-TextValue $wsOut "D2" "00044002 - Al Piemonte Ford Sales Inc"
Set-TextValue $wsOut "E2" "`$33,044.00"
Set-TextValue $wsOut "F2" "PON00212643/0"
Set-TextValue $wsOut "G2" "INV00212643"

# Drop the stale PO Amount_posted data cell and the two now-reconciled rows
$wsOut.Rows("3:4").Delete()
$wsOut.Columns(8).ClearFormats()
$wsOut.Range("H2").Value = ""

# Shrink the column widths to fit the new (shorter) content
$wsOut.Columns(1).ColumnWidth = 4.84375
$wsOut.Columns(2).ColumnWidth = 9.98046875
$wsOut.Columns(3).ColumnWidth = 4.2578125
$wsOut.Columns(4).ColumnWidth = 7.6171875
$wsOut.Columns(6).ColumnWidth = 5.05078125
$wsOut.Columns(7).ColumnWidth = 12.43359375

# Keep the dependent lookup sheets in sync with the new unit number
Set-TextValue $wsUnit "A2" "00997989"
Set-TextValue $wsPO "A2" "00997989"

# --- Active selections / active sheet ---
$wsLogin.Range("A18").Select()
$wsUnit.Range("A2").Select()
$wsPO.Range("A3").Select()
$wsAccept.Range("I19").Select()
$wsOut.Range("A1").Select()

# PO_Detail is the sheet left active/selected after the reconcile pass
$wsPO.Select()
